$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely, shifting row 3 up into row 2's position
$ws.Rows(2).Delete()
